$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 127.14286
$ws.Range("I8").Value = 127.14286
$ws.Range("K8").Value = 381.42858
$ws.Range("M8").Value = -242.42858
# Row 9
$ws.Range("H9").Value = 64.454544
$ws.Range("I9").Value = 45.444443
$ws.Range("K9").Value = 45.444443
$ws.Range("M9").Value = 123.555557
# Row 19
$ws.Range("H19").Value = 5761.1924
$ws.Range("I19").Value = 2117
$ws.Range("J19").Value = 8433.6
$ws.Range("K19").Value = 2117
$ws.Range("L19").Value = 8433.6
$ws.Range("M19").Value = -1942
$ws.Range("N19").Value = -8783.6
# Row 33
$ws.Range("H33").Value = 19609488
$ws.Range("I33").Value = 22223786
$ws.Range("J33").Value = 2249.5
$ws.Range("K33").Value = 22223786
$ws.Range("L33").Value = 2249.5
$ws.Range("M33").Value = -22223557
$ws.Range("N33").Value = -2707.5
# Row 112
$ws.Range("H112").Value = 2147.3044
$ws.Range("J112").Value = 2470.1765
$ws.Range("L112").Value = 7410.529500000001
$ws.Range("N112").Value = -9626.529500000001
# Row 116
$ws.Range("H116").Value = 6182.5
$ws.Range("I116").Value = 6275.933
$ws.Range("J116").Value = 6055.091
$ws.Range("K116").Value = 6275.933
$ws.Range("L116").Value = 6055.091
$ws.Range("M116").Value = -2833.933
$ws.Range("N116").Value = -12939.091

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2595.28
$ws.Range("I2").Value = 1995.826
$ws.Range("K2").Value = 1995.826
$ws.Range("M2").Value = -1882.826
# Row 4
$ws.Range("H4").Value = 2909.24
$ws.Range("I4").Value = 2488.3914
$ws.Range("J4").Value = 7749
$ws.Range("K4").Value = 2488.3914
$ws.Range("L4").Value = 7749
$ws.Range("M4").Value = -2372.3914
$ws.Range("N4").Value = -7981
# Row 5
$ws.Range("H5").Value = 323.75
$ws.Range("I5").Value = 398.33334
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 398.33334
$ws.Range("L5").Value = 100
$ws.Range("M5").Value = -286.33334
$ws.Range("N5").Value = -324
# Row 32
$ws.Range("H32").Value = 10490.44
$ws.Range("I32").Value = 7443.7144
$ws.Range("J32").Value = 20690.348
$ws.Range("K32").Value = 7443.7144
$ws.Range("L32").Value = 20690.348
$ws.Range("M32").Value = -7156.7144
$ws.Range("N32").Value = -21264.348
# Row 37
$ws.Range("H37").Value = 39999
$ws.Range("J37").Value = 59999
$ws.Range("L37").Value = 59999
$ws.Range("N37").Value = -60545
# Row 45
$ws.Range("H45").Value = 2448.0208
$ws.Range("I45").Value = 1395.5652
$ws.Range("K45").Value = 1395.5652
$ws.Range("M45").Value = -1018.5652
# Row 63
$ws.Range("H63").Value = 7781.25
$ws.Range("I63").Value = 7062.5
$ws.Range("K63").Value = 7062.5
$ws.Range("M63").Value = -6376.5
# Row 66
$ws.Range("H66").Value = 7781.25
$ws.Range("I66").Value = 7062.5
$ws.Range("K66").Value = 35312.5
$ws.Range("M66").Value = -31880.5
# Row 116
$ws.Range("H116").Value = 2595.28
$ws.Range("I116").Value = 1995.826
$ws.Range("K116").Value = 1995.826
$ws.Range("M116").Value = 298.174
# Row 122
$ws.Range("H122").Value = 6641.7
$ws.Range("I122").Value = 5791
$ws.Range("J122").Value = 6854.375
$ws.Range("K122").Value = 17373
$ws.Range("L122").Value = 20563.125
$ws.Range("M122").Value = -14923
$ws.Range("N122").Value = -25463.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2595.28
$ws.Range("I3").Value = 1995.826
$ws.Range("K3").Value = 1995.826
$ws.Range("M3").Value = -1881.826
# Row 4
$ws.Range("H4").Value = 323.75
$ws.Range("I4").Value = 398.33334
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 398.33334
$ws.Range("L4").Value = 100
$ws.Range("M4").Value = -283.33334
$ws.Range("N4").Value = -330
# Row 11
$ws.Range("H11").Value = 1950
$ws.Range("I11").Value = 2250
$ws.Range("J11").Value = 1500
$ws.Range("K11").Value = 2250
$ws.Range("L11").Value = 1500
$ws.Range("M11").Value = -2110
$ws.Range("N11").Value = -1780
# Row 20
$ws.Range("H20").Value = 3616.3
$ws.Range("I20").Value = 3494.6924
$ws.Range("J20").Value = 3842.1428
$ws.Range("K20").Value = 3494.6924
$ws.Range("L20").Value = 3842.1428
$ws.Range("M20").Value = -3247.6924
$ws.Range("N20").Value = -4336.1428
# Row 105
$ws.Range("H105").Value = 3150.1428
$ws.Range("I105").Value = 3079.5
$ws.Range("J105").Value = 3326.75
$ws.Range("K105").Value = 3079.5
$ws.Range("L105").Value = 3326.75
$ws.Range("M105").Value = -1332.5
$ws.Range("N105").Value = -6820.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 678
$ws.Range("I16").Value = 678
$ws.Range("K16").Value = 678
$ws.Range("M16").Value = -391
# Row 22
$ws.Range("H22").Value = 499.2353
$ws.Range("I22").Value = 273.36365
$ws.Range("J22").Value = 913.3333
$ws.Range("K22").Value = 273.36365
$ws.Range("L22").Value = 913.3333
$ws.Range("M22").Value = 76.63634999999999
$ws.Range("N22").Value = -1613.3333
# Row 62
$ws.Range("H62").Value = 8779.5
$ws.Range("J62").Value = 8113.4443
$ws.Range("L62").Value = 8113.4443
$ws.Range("N62").Value = -9361.444299999999
# Row 65
$ws.Range("H65").Value = 8779.5
$ws.Range("J65").Value = 8113.4443
$ws.Range("L65").Value = 40567.2215
$ws.Range("N65").Value = -46807.2215
# Row 113
$ws.Range("H113").Value = 678
$ws.Range("I113").Value = 678
$ws.Range("K113").Value = 678
$ws.Range("M113").Value = 1492

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 23281664
$ws.Range("I4").Value = 27806306
$ws.Range("K4").Value = 83418918
$ws.Range("M4").Value = -83418806
# Row 69
$ws.Range("H69").Value = 3003
$ws.Range("I69").Value = 2012
$ws.Range("J69").Value = 3333.3333
$ws.Range("K69").Value = 6036
$ws.Range("L69").Value = 9999.999899999999
$ws.Range("M69").Value = -5225
$ws.Range("N69").Value = -11621.9999
# Row 72
$ws.Range("H72").Value = 3003
$ws.Range("I72").Value = 2012
$ws.Range("J72").Value = 3333.3333
$ws.Range("K72").Value = 18108
$ws.Range("L72").Value = 29999.9997
$ws.Range("M72").Value = -14052
$ws.Range("N72").Value = -38111.9997
# Row 107
$ws.Range("H107").Value = 1475.5333
$ws.Range("J107").Value = 1671.32
$ws.Range("L107").Value = 5013.96
$ws.Range("N107").Value = -8853.959999999999
# Row 113
$ws.Range("H113").Value = 184739.8
$ws.Range("I113").Value = 20000
$ws.Range("K113").Value = 60000
$ws.Range("M113").Value = -57830
# Row 129
$ws.Range("H129").Value = 1326.2858
$ws.Range("I129").Value = 853.44446
$ws.Range("J129").Value = 2177.4
$ws.Range("K129").Value = 2560.33338
$ws.Range("L129").Value = 6532.200000000001
$ws.Range("M129").Value = 2439.66662
$ws.Range("N129").Value = -16532.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 5542.7617
$ws.Range("I122").Value = 4897.875
$ws.Range("K122").Value = 14693.625
$ws.Range("M122").Value = -12243.625
# Row 132
$ws.Range("H132").Value = 9112
$ws.Range("I132").Value = 10670
$ws.Range("J132").Value = 8333
$ws.Range("K132").Value = 32010
$ws.Range("L132").Value = 24999
$ws.Range("M132").Value = -29480
$ws.Range("N132").Value = -30059

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 6206.231
$ws.Range("I132").Value = 3801
$ws.Range("J132").Value = 9012.333000000001
$ws.Range("K132").Value = 11403
$ws.Range("L132").Value = 27036.999
$ws.Range("M132").Value = -8873
$ws.Range("N132").Value = -32096.999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 17790.375
$ws.Range("I52").Value = 17790.375
$ws.Range("K52").Value = 17790.375
$ws.Range("M52").Value = -17564.375
# Row 124
$ws.Range("H124").Value = 69000
$ws.Range("J124").Value = 69000
$ws.Range("L124").Value = 69000
$ws.Range("N124").Value = -78820
# Row 126
$ws.Range("H126").Value = 27073.389
$ws.Range("I126").Value = 36360.082
$ws.Range("K126").Value = 109080.246
$ws.Range("M126").Value = -106610.246
